$wb = $excel.ActiveWorkbook

# --- 1. Clear B4, B7, B8, B9 on "ODI Batting" sheet ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B4").ClearContents()
$batting.Range("B7").ClearContents()
$batting.Range("B8").ClearContents()
$batting.Range("B9").ClearContents()

# --- 2. Add new worksheet "ODI Batting Extra" after "ODI Bowling" ---
$bowling = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowling)
$extra.Name = "ODI Batting Extra"

# Header row (text labels)
$extra.Cells.Item(1, 1).NumberFormat = "@"
$extra.Cells.Item(1, 1).Value = "MATCH_CODE"
$extra.Cells.Item(1, 2).NumberFormat = "@"
$extra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$extra.Cells.Item(1, 3).NumberFormat = "@"
$extra.Cells.Item(1, 3).Value = "NUM_4"
$extra.Cells.Item(1, 4).NumberFormat = "@"
$extra.Cells.Item(1, 4).Value = "NUM_6"
$extra.Cells.Item(1, 5).NumberFormat = "@"
$extra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Cells.Item(1, 6).NumberFormat = "@"
$extra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# Data rows. Columns A, C, D, E hold digit-like text ("4628", "0", "0.41%") that must
# stay text rather than being auto-coerced to numbers, so force NumberFormat="@" first.
# Column B (BATTING_POSITION) holds real numbers, and F holds plain text ("NO").
$data = @(
    @("4628", 11,   "0", "0", "",      "NO"),
    @("4679", 9,    "0", "0", "",      "NO"),
    @("4682", 10,   "",  "",  "",      "NO"),
    @("4685", "",   "",  "",  "",      "NO"),
    @("4717", 10,   "0", "0", "0.41%", "NO"),
    @("4726", "",   "",  "",  "",      "NO"),
    @("4729", "",   "",  "",  "",      "NO"),
    @("4734", 9,    "",  "",  "",      "NO")
)

$rowIndex = 2
foreach ($row in $data) {
    $extra.Cells.Item($rowIndex, 1).NumberFormat = "@"
    $extra.Cells.Item($rowIndex, 1).Value = $row[0]

    if ($row[1] -eq "") {
        # Touch the cell first so an empty assignment still leaves a (blank) cell
        # behind instead of no cell at all.
        $extra.Cells.Item($rowIndex, 2).NumberFormat = "General"
    }
    $extra.Cells.Item($rowIndex, 2).Value = $row[1]

    $extra.Cells.Item($rowIndex, 3).NumberFormat = "@"
    $extra.Cells.Item($rowIndex, 3).Value = $row[2]

    $extra.Cells.Item($rowIndex, 4).NumberFormat = "@"
    $extra.Cells.Item($rowIndex, 4).Value = $row[3]

    $extra.Cells.Item($rowIndex, 5).NumberFormat = "@"
    $extra.Cells.Item($rowIndex, 5).Value = $row[4]

    $extra.Cells.Item($rowIndex, 6).Value = $row[5]

    $rowIndex++
}

$wb.Save()
